$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the hint/info cell in D4 --------------------------------------
# Was: "Test-Text:info;Hinweis-AMA Daten werden von der Finanz verwertet, ..."
# Now: "HINWEIS --> AMA Daten werden von der Finanz verwertet, ..."
$ws.Range("D4").Value = "HINWEIS --> AMA Daten werden von der Finanz verwertet, saubere Erstdaten helfen bei Plausibilitätsprüfung:info;Flächenaufstellung:select(Eigen-,Pacht,Mitbewirtschaftung):pflicht"

# --- Add "Upload:checkbox:pflicht" to D9:D12 ------------------------------
# These rows (A2.1 Fragen zum Betrieb sub-questions) previously had no
# value in column D; give each an upload/checkbox requirement, using the
# same text number-format ("@") already applied to the rest of column D.
foreach ($row in 9..12) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.Value = "Upload:checkbox:pflicht"
    $cell.NumberFormat = "@"
}

# --- Update the saved view: scroll position + active selection -----------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("D8").Select()
